$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a mis-entered value in row 7 (C7 was 85.747, should be 82.747)
$ws.Range("C7").Value = 82.747

# Add a new row (row 8) that computes the weekly "new cases" figures
# by taking the week-over-week delta (multiplied by 1000) of row 7.
$ws.Range("A8").Value = "新增"
$ws.Range("B8").Value = 0

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $cur = $cols[$i]
    if ($i -eq 0) {
        $prev = "B"
    } else {
        $prev = $cols[$i - 1]
    }
    $ws.Range($cur + "8").Formula = "=" + $cur + "7*1000-" + $prev + "7*1000"
}

# Leave the selection where the new data was entered
$ws.Range("O21").Select()
